$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 115 (pushing the existing
# rows 115-163 down to 116-164, and growing the used range to A1:R164).
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row with its own data.
$ws.Cells.Item(115, 1).Value = 7
$ws.Cells.Item(115, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115, 3).Value = "Ñuble"
$ws.Cells.Item(115, 4).Value = 44553
$ws.Cells.Item(115, 5).Value = 16
$ws.Cells.Item(115, 6).Value = 100112032
$ws.Cells.Item(115, 7).Value = "Zapallo italiano"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 60
$ws.Cells.Item(115, 11).Value = 4000
$ws.Cells.Item(115, 12).Value = 4500
$ws.Cells.Item(115, 13).Value = 4250
$ws.Cells.Item(115, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(115, 15).Value = "Región del Maule"
$ws.Cells.Item(115, 16).Value = 71
$ws.Cells.Item(115, 17).Value = 60
$ws.Cells.Item(115, 18).Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(115, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat
